$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 37037500
$ws.Range("I33").Value = 138.85
$ws.Range("K33").Value = 138.85
$ws.Range("M33").Value = 90.15000000000001

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2998.0938
$ws.Range("I64").Value = 2738.75
$ws.Range("J64").Value = 3084.5417
$ws.Range("K64").Value = 2738.75
$ws.Range("L64").Value = 3084.5417
$ws.Range("M64").Value = -2490.75
$ws.Range("N64").Value = -3580.5417

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 2998.0938
$ws.Range("I67").Value = 2738.75
$ws.Range("J67").Value = 3084.5417
$ws.Range("K67").Value = 2738.75
$ws.Range("L67").Value = 3084.5417
$ws.Range("M67").Value = -1880.75
$ws.Range("N67").Value = -4800.5417

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1208965.2
$ws.Range("I137").Value = 1281.0883
$ws.Range("J137").Value = 4630737
$ws.Range("K137").Value = 3843.2649
$ws.Range("L137").Value = 13892211
$ws.Range("M137").Value = -1293.2649
$ws.Range("N137").Value = -13897311

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 895.86957
$ws.Range("I141").Value = 804.6667
$ws.Range("K141").Value = 2414.0001
$ws.Range("M141").Value = 2765.9999

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2000
$ws.Range("I2").Value = 1861.1111
$ws.Range("J2").Value = 2357.1428
$ws.Range("K2").Value = 1861.1111
$ws.Range("L2").Value = 2357.1428
$ws.Range("M2").Value = -1748.1111
$ws.Range("N2").Value = -2583.1428

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2274.82
$ws.Range("I32").Value = 2051.3813
$ws.Range("J32").Value = 9499.333000000001
$ws.Range("K32").Value = 2051.3813
$ws.Range("L32").Value = 9499.333000000001
$ws.Range("M32").Value = -1764.3813
$ws.Range("N32").Value = -10073.333

# ARM row 52
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 23446.666
$ws.Range("J52").Value = 23446.666
$ws.Range("L52").Value = 23446.666
$ws.Range("N52").Value = -24082.666

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 38526.3
$ws.Range("I74").Value = 60263.707
$ws.Range("J74").Value = 10100.462
$ws.Range("K74").Value = 60263.707
$ws.Range("L74").Value = 10100.462
$ws.Range("M74").Value = -59389.707
$ws.Range("N74").Value = -11848.462

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 38526.3
$ws.Range("I77").Value = 60263.707
$ws.Range("J77").Value = 10100.462
$ws.Range("K77").Value = 301318.535
$ws.Range("L77").Value = 50502.31
$ws.Range("M77").Value = -296950.535
$ws.Range("N77").Value = -59238.31

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2000
$ws.Range("I116").Value = 1861.1111
$ws.Range("J116").Value = 2357.1428
$ws.Range("K116").Value = 1861.1111
$ws.Range("L116").Value = 2357.1428
$ws.Range("M116").Value = 432.8888999999999
$ws.Range("N116").Value = -6945.1428

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1589.6086
$ws.Range("I122").Value = 1284.75
$ws.Range("J122").Value = 1922.1818
$ws.Range("K122").Value = 3854.25
$ws.Range("L122").Value = 5766.5454
$ws.Range("M122").Value = -1404.25
$ws.Range("N122").Value = -10666.5454

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1995.0333
$ws.Range("I132").Value = 1819.2858
$ws.Range("J132").Value = 2405.111
$ws.Range("K132").Value = 5457.857400000001
$ws.Range("L132").Value = 7215.333
$ws.Range("M132").Value = -2927.857400000001
$ws.Range("N132").Value = -12275.333

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2000
$ws.Range("I3").Value = 1861.1111
$ws.Range("J3").Value = 2357.1428
$ws.Range("K3").Value = 1861.1111
$ws.Range("L3").Value = 2357.1428
$ws.Range("M3").Value = -1747.1111
$ws.Range("N3").Value = -2585.1428

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 583271.4
$ws.Range("I134").Value = 933872.9
$ws.Range("J134").Value = 3430.5
$ws.Range("K134").Value = 2801618.7
$ws.Range("L134").Value = 10291.5
$ws.Range("M134").Value = -2799083.7
$ws.Range("N134").Value = -15361.5

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 670.2105
$ws.Range("J22").Value = 1003.1818
$ws.Range("L22").Value = 1003.1818
$ws.Range("N22").Value = -1703.1818

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9260441
$ws.Range("I31").Value = 1048.8445
$ws.Range("J31").Value = 55557404
$ws.Range("K31").Value = 1048.8445
$ws.Range("L31").Value = 55557404
$ws.Range("M31").Value = -753.8444999999999
$ws.Range("N31").Value = -55557994

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 9260441
$ws.Range("I34").Value = 1048.8445
$ws.Range("J34").Value = 55557404
$ws.Range("K34").Value = 1048.8445
$ws.Range("L34").Value = 55557404
$ws.Range("M34").Value = -846.8444999999999
$ws.Range("N34").Value = -55557808

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 807074.9
$ws.Range("I132").Value = 1983.45
$ws.Range("J132").Value = 6174351
$ws.Range("K132").Value = 5950.35
$ws.Range("L132").Value = 18523053
$ws.Range("M132").Value = -3420.35
$ws.Range("N132").Value = -18528113

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2085.8594
$ws.Range("I134").Value = 2342.7234
$ws.Range("J134").Value = 1375.7059
$ws.Range("K134").Value = 7028.1702
$ws.Range("L134").Value = 4127.1177
$ws.Range("M134").Value = -4493.1702
$ws.Range("N134").Value = -9197.117699999999

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 718.9545000000001
$ws.Range("I122").Value = 589.2353000000001
$ws.Range("J122").Value = 1160
$ws.Range("K122").Value = 5303.117700000001
$ws.Range("L122").Value = 10440
$ws.Range("M122").Value = -2853.117700000001
$ws.Range("N122").Value = -15340

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1919.5714
$ws.Range("I102").Value = 1737.1666
$ws.Range("K102").Value = 1737.1666
$ws.Range("M102").Value = -115.1666

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1686.4546
$ws.Range("I16").Value = 1475.1
$ws.Range("J16").Value = 3800
$ws.Range("K16").Value = 1475.1
$ws.Range("L16").Value = 3800
$ws.Range("M16").Value = -1305.1
$ws.Range("N16").Value = -4140

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2630.8
$ws.Range("I40").Value = 2501.3333
$ws.Range("J40").Value = 2825
$ws.Range("K40").Value = 2501.3333
$ws.Range("L40").Value = 2825
$ws.Range("M40").Value = -2365.3333
$ws.Range("N40").Value = -3097

# LTW row 51
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

# LTW row 127
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 34221.6
$ws.Range("J127").Value = 34221.6
$ws.Range("L127").Value = 34221.6
$ws.Range("N127").Value = -44141.6

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4576.595
$ws.Range("I132").Value = 4577.0347
$ws.Range("J132").Value = 4575.615
$ws.Range("K132").Value = 13731.1041
$ws.Range("L132").Value = 13726.845
$ws.Range("M132").Value = -11201.1041
$ws.Range("N132").Value = -18786.845

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1685.1915
$ws.Range("I136").Value = 1048.4
$ws.Range("K136").Value = 3145.2
$ws.Range("M136").Value = -595.2000000000003

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2032.9592
$ws.Range("I132").Value = 2059.5588
$ws.Range("J132").Value = 1972.6666
$ws.Range("K132").Value = 6178.676399999999
$ws.Range("L132").Value = 5917.9998
$ws.Range("M132").Value = -3648.676399999999
$ws.Range("N132").Value = -10977.9998
